# Add the "rest" position row back into the RobotPositions sheet (row 14):
#  - A14 becomes a label cell ("rest") formatted like the other position
#    labels in column A (A4:A13)
#  - C14/D14/E14 get their numeric values back (500 / 150 / 340)
#
# The workbook is already open.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give A14 the same cell formatting as the other row labels (A13 is the
# closest fully-styled label cell) before writing the text, so the label
# style (border/fill/text-format) matches the rest of the table.
$ws.Range("A13").Copy()
$ws.Range("A14").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A14").Value = "rest"

# Restore the numeric values for this position's s1 (C), s2 (D) and s3 (E)
# joint coordinates.
$ws.Range("C14").Value = 500
$ws.Range("D14").Value = 150
$ws.Range("E14").Value = 340
